# StatBudget.docx - "import rappel fait from xlsx"
#
# The document contains a single table tracking the monthly budget
# consumption ("PAIEMENT MENSUELLE"), the back-pay / recall amounts
# ("RAPPEL"), the combined totals ("TOTAL") and the running balance
# ("SOLD"), one row per month (JAN..DEC).
#
# Table columns (1-based, as seen by Word's Table.Cell(row, col)):
#   1 = MONTANT ALLOUE (allocated amount, only filled on the JAN row)
#   2 = MOIS (month label - unchanged)
#   3 = PAIEMENT MENSUELLE / NOMBRES
#   4 = PAIEMENT MENSUELLE / MONTANT CONSOMME
#   5 = RAPPEL / NOMBRES
#   6 = RAPPEL / MONTANT CONSOMME
#   7 = TOTAL / NIMBRES TOTAL
#   8 = TOTAL / MONTANT CONSOMME TOTAL
#   9 = SOLD
#
# Table rows: row 1-2 = headers, rows 3-14 = JAN..DEC data rows.
#
# This refreshes the whole table with numbers re-imported from the
# source spreadsheet: JAN/FEV are adjusted, MARS/AVRIL now carry a
# RAPPEL amount, and MAI..DEC are zeroed out (no consumption recorded
# for those months in the refreshed import), with the SOLD simply
# carrying forward the balance after AVRIL for the remaining months.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-Cell($row, $col, $value) {
    $t.Cell($row, $col).Range.Text = $value
}

# Row 3 - JAN
Set-Cell 3 1 "52 769 000,00"
Set-Cell 3 3 "4437"
Set-Cell 3 4 "3 993 300,00"
Set-Cell 3 7 "4437"
Set-Cell 3 8 "3 993 300,00"
Set-Cell 3 9 "48 775 700,00"

# Row 4 - FEV
Set-Cell 4 3 "3144"
Set-Cell 4 4 "2 829 600,00"
Set-Cell 4 7 "3144"
Set-Cell 4 8 "2 829 600,00"
Set-Cell 4 9 "45 946 100,00"

# Row 5 - MARS (gains a RAPPEL amount)
Set-Cell 5 3 "3868"
Set-Cell 5 4 "3 481 200,00"
Set-Cell 5 5 "724"
Set-Cell 5 6 "651 600,00"
Set-Cell 5 7 "4592"
Set-Cell 5 8 "4 132 800,00"
Set-Cell 5 9 "41 813 300,00"

# Row 6 - AVRIL (gains a RAPPEL amount)
Set-Cell 6 3 "4053"
Set-Cell 6 4 "3 647 700,00"
Set-Cell 6 5 "495"
Set-Cell 6 6 "445 500,00"
Set-Cell 6 7 "4548"
Set-Cell 6 8 "4 093 200,00"
Set-Cell 6 9 "37 720 100,00"

# Row 7 - MAI (zeroed out, SOLD carries forward)
Set-Cell 7 3 "0"
Set-Cell 7 4 "0,00"
Set-Cell 7 7 "0"
Set-Cell 7 8 "0,00"
Set-Cell 7 9 "37 720 100,00"

# Row 8 - JUIN (zeroed out, SOLD carries forward)
Set-Cell 8 3 "0"
Set-Cell 8 4 "0,00"
Set-Cell 8 7 "0"
Set-Cell 8 8 "0,00"
Set-Cell 8 9 "37 720 100,00"

# Row 9 - JUIL (zeroed out, SOLD carries forward)
Set-Cell 9 3 "0"
Set-Cell 9 4 "0,00"
Set-Cell 9 5 "0"
Set-Cell 9 7 "0"
Set-Cell 9 8 "0,00"
Set-Cell 9 9 "37 720 100,00"

# Row 10 - AOUT (zeroed out, SOLD carries forward)
Set-Cell 10 3 "0"
Set-Cell 10 4 "0,00"
Set-Cell 10 7 "0"
Set-Cell 10 8 "0,00"
Set-Cell 10 9 "37 720 100,00"

# Row 11 - SEPT (zeroed out, SOLD carries forward)
Set-Cell 11 3 "0"
Set-Cell 11 4 "0,00"
Set-Cell 11 5 "0"
Set-Cell 11 6 "0,00"
Set-Cell 11 7 "0"
Set-Cell 11 8 "0,00"
Set-Cell 11 9 "37 720 100,00"

# Row 12 - OCT (zeroed out, SOLD carries forward)
Set-Cell 12 3 "0"
Set-Cell 12 4 "0,00"
Set-Cell 12 5 "0"
Set-Cell 12 6 "0,00"
Set-Cell 12 7 "0"
Set-Cell 12 8 "0,00"
Set-Cell 12 9 "37 720 100,00"

# Row 13 - NOV (zeroed out, SOLD carries forward)
Set-Cell 13 3 "0"
Set-Cell 13 4 "0,00"
Set-Cell 13 5 "0"
Set-Cell 13 6 "0,00"
Set-Cell 13 7 "0"
Set-Cell 13 8 "0,00"
Set-Cell 13 9 "37 720 100,00"

# Row 14 - DEC (zeroed out, SOLD carries forward)
Set-Cell 14 3 "0"
Set-Cell 14 4 "0,00"
Set-Cell 14 7 "0"
Set-Cell 14 8 "0,00"
Set-Cell 14 9 "37 720 100,00"
